# Add a new "Comparison" layout slide (Data Analysis: Pros / Cons) as the
# new last slide (slide 11) of the deck.

$p = $ppt.ActivePresentation

# --- Step 1: create the slide with a blank layout first, then add/remove a
# couple of dummy shapes. This nudges the per-slide shape-id counter forward
# by 2 so that, once we switch to the real "Comparison" layout, the 5
# placeholders it materializes land on ids 4-8 (matching the authored file)
# instead of 2-6.
$layoutBlank = $p.SlideMaster.CustomLayouts.Item(7)
$slide = $p.Slides.AddSlide($p.Slides.Count + 1, $layoutBlank)

$dummy1 = $slide.Shapes.AddTextbox(1, 0, 0, 10, 10)
$dummy2 = $slide.Shapes.AddTextbox(1, 0, 0, 10, 10)
$dummy1.Delete()
$dummy2.Delete()

# --- Step 2: switch the slide to the real "Comparison" layout (5th custom
# layout on the master: Title, Text/Content halves, Text/Content quarters).
$layoutComparison = $p.SlideMaster.CustomLayouts.Item(5)
$slide.CustomLayout = $layoutComparison

# --- Step 3: name + fill in each placeholder to match the authored shapes.
$shTitle = $slide.Shapes.Item(1)
$shTitle.Name = "Title 3"
$shTitle.TextFrame.TextRange.Text = "Data Analysis"

$shProsLabel = $slide.Shapes.Item(2)
$shProsLabel.Name = "Text Placeholder 4"
$shProsLabel.TextFrame.TextRange.Text = "Pros"

$shProsBody = $slide.Shapes.Item(3)
$shProsBody.Name = "Content Placeholder 5"
$prTr = $shProsBody.TextFrame.TextRange
$prTr.Text = "Improved decision-making"
$prTr.LanguageID = "en-GB"
$prosRest = @("Cost-effective", "Innovation", "Competitive advantage", "Personalization")
foreach ($line in $prosRest) {
    $added = $shProsBody.TextFrame.TextRange.InsertAfter("`r" + $line)
    $added.LanguageID = "en-GB"
}

$shConsLabel = $slide.Shapes.Item(4)
$shConsLabel.Name = "Text Placeholder 6"
$shConsLabel.TextFrame.TextRange.Text = "Cons"

$shConsBody = $slide.Shapes.Item(5)
$shConsBody.Name = "Content Placeholder 7"
$coTr = $shConsBody.TextFrame.TextRange
$coTr.Text = "Data quality"
$coTr.LanguageID = "en-GB"
$consRest = @("Privacy concerns", "Complexity", "Bias", "Interpretation")
foreach ($line in $consRest) {
    $added = $shConsBody.TextFrame.TextRange.InsertAfter("`r" + $line)
    $added.LanguageID = "en-GB"
}
